$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before row 154 (shifts existing rows 154-227 down to 156-229)
$ws.Rows("154:155").Insert()

# Populate the two new rows with the new weekly price data (fecha 2022-02-17)
$ws.Cells.Item(154,1).Value = 2
$ws.Cells.Item(154,2).Value = "Comercializadora del Agro de Limarí"
$ws.Cells.Item(154,3).Value = "Coquimbo"
$ws.Cells.Item(154,4).Value = 44609
$ws.Cells.Item(154,5).Value = 4
$ws.Cells.Item(154,6).Value = 100112021
$ws.Cells.Item(154,7).Value = "Ají"
$ws.Cells.Item(154,8).Value = "Americana (o)"
$ws.Cells.Item(154,9).Value = "Primera"
$ws.Cells.Item(154,10).Value = 300
$ws.Cells.Item(154,11).Value = 9000
$ws.Cells.Item(154,12).Value = 11000
$ws.Cells.Item(154,13).Value = 10000
$ws.Cells.Item(154,14).Value = "`$/caja 25 kilos"
$ws.Cells.Item(154,15).Value = "Provincia de Limarí"
$ws.Cells.Item(154,16).Value = 400
$ws.Cells.Item(154,17).Value = 25
$ws.Cells.Item(154,18).Value = "Hortaliza"
$ws.Cells.Item(155,1).Value = 2
$ws.Cells.Item(155,2).Value = "Comercializadora del Agro de Limarí"
$ws.Cells.Item(155,3).Value = "Coquimbo"
$ws.Cells.Item(155,4).Value = 44609
$ws.Cells.Item(155,5).Value = 4
$ws.Cells.Item(155,6).Value = 100112021
$ws.Cells.Item(155,7).Value = "Ají"
$ws.Cells.Item(155,8).Value = "Inferno"
$ws.Cells.Item(155,9).Value = "Primera"
$ws.Cells.Item(155,10).Value = 160
$ws.Cells.Item(155,11).Value = 18000
$ws.Cells.Item(155,12).Value = 20000
$ws.Cells.Item(155,13).Value = 19000
$ws.Cells.Item(155,14).Value = "`$/caja 25 kilos"
$ws.Cells.Item(155,15).Value = "Provincia de Limarí"
$ws.Cells.Item(155,16).Value = 760
$ws.Cells.Item(155,17).Value = 25
$ws.Cells.Item(155,18).Value = "Hortaliza"
